$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = -0.010862383516112573
$ws.Range("B1").Value = 0.010724236323085279
$ws.Range("A2").Value = 0.03538265726256462
$ws.Range("B2").Value = -0.03602170191154919
$ws.Range("A3").Value = 0.14456276755547748
$ws.Range("B3").Value = -0.14525711302865218
$ws.Range("A4").Value = -0.13872635237535746
$ws.Range("B4").Value = 0.13808258844101573
$ws.Range("A5").Value = -0.13208258929532413
$ws.Range("B5").Value = 0.1308026531361426
$ws.Range("A6").Value = -0.029977692536763634
$ws.Range("B6").Value = 0.02996464408558852
$ws.Range("A7").Value = -0.00996464510875228
$ws.Range("B7").Value = 0.009935275865094795
$ws.Range("A8").Value = 0.010064723110573048
$ws.Range("B8").Value = -0.010148116444994137
$ws.Range("A9").Value = 0.01614811556747675
$ws.Range("B9").Value = -0.016252587645658068
$ws.Range("A10").Value = -0.015049004247778441
$ws.Range("B10").Value = 0.01504977678158781
$ws.Range("A11").Value = -0.010549777640804336
$ws.Range("B11").Value = 0.010548646909114723
$ws.Range("A12").Value = -0.004548647784290871
$ws.Range("B12").Value = 0.004541954658476843
$ws.Range("A13").Value = 0.001458044465355357
$ws.Range("B13").Value = -0.0014584960702306304
$ws.Range("A14").Value = 0.013458495131996706
$ws.Range("B14").Value = -0.013465008365534104
$ws.Range("A15").Value = 0.01946500749059421
$ws.Range("B15").Value = -0.019483415053097808
$ws.Range("A16").Value = 0.025483414180440978
$ws.Range("B16").Value = -0.0255412987040069
$ws.Range("A17").Value = -0.009003690709813661
$ws.Range("B17").Value = 0.008999999102116796
$ws.Range("A18").Value = -0.03610689448284532
$ws.Range("B18").Value = 0.036095891612639264
$ws.Range("A19").Value = -0.027095892499889196
$ws.Range("B19").Value = 0.027013039798366556
$ws.Range("A20").Value = -0.018013040693102056
$ws.Range("B20").Value = 0.018004197303698177
$ws.Range("A21").Value = -0.00900419819947551
$ws.Range("B21").Value = 0.008999999103496137
$ws.Range("A22").Value = -0.09954106652859451
$ws.Range("B22").Value = 0.09922323497589858
$ws.Range("A23").Value = -0.060855404193408624
$ws.Range("B23").Value = 0.060638556761275275
$ws.Range("A24").Value = -0.04212338418599604
$ws.Range("B24").Value = 0.04199999874890192
$ws.Range("A25").Value = -0.09482496147653308
$ws.Range("B25").Value = 0.094591968789139
$ws.Range("A26").Value = -0.08859196966910687
$ws.Range("B26").Value = 0.08828935867786569
$ws.Range("A27").Value = -0.08228935956206573
$ws.Range("B27").Value = 0.08124547830333828
$ws.Range("A28").Value = -0.07524547920727365
$ws.Range("B28").Value = 0.07452244621002624
$ws.Range("A29").Value = -0.06252244718967148
$ws.Range("B29").Value = 0.06216504933545508
$ws.Range("A30").Value = -0.04216505040411134
$ws.Range("B30").Value = 0.042018237447402296
$ws.Range("A31").Value = -0.027018238473052847
$ws.Range("B31").Value = 0.027000501676461397
$ws.Range("A32").Value = -0.0060005027655103405
$ws.Range("B32").Value = 0.005999999066200523
